$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4060.5
$ws.Range("I113").Value = 2853.75
$ws.Range("J113").Value = 4865
$ws.Range("K113").Value = 2853.75
$ws.Range("L113").Value = 4865
$ws.Range("M113").Value = 400.25
$ws.Range("N113").Value = -11373

$ws.Range("H125").Value = 83333550
$ws.Range("I125").Value = 400
$ws.Range("J125").Value = 125000130
$ws.Range("K125").Value = 3600
$ws.Range("L125").Value = 1125001170
$ws.Range("M125").Value = -1140
$ws.Range("N125").Value = -1125006090

$ws.Range("H132").Value = 3050.2727
$ws.Range("I132").Value = 1884.5641
$ws.Range("J132").Value = 5891.6875
$ws.Range("K132").Value = 5653.692300000001
$ws.Range("L132").Value = 17675.0625
$ws.Range("M132").Value = -3123.692300000001
$ws.Range("N132").Value = -22735.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1018.74
$ws.Range("I32").Value = 1008.8283
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 1008.8283
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -721.8283
$ws.Range("N32").Value = -2574

$ws.Range("H61").Value = 249098.11
$ws.Range("I61").Value = 212728.28
$ws.Range("J61").Value = 306408.12
$ws.Range("K61").Value = 212728.28
$ws.Range("L61").Value = 306408.12
$ws.Range("M61").Value = -212516.28
$ws.Range("N61").Value = -306832.12

$ws.Range("H63").Value = 7986.846
$ws.Range("I63").Value = 9172.4
$ws.Range("J63").Value = 4035
$ws.Range("K63").Value = 9172.4
$ws.Range("L63").Value = 4035
$ws.Range("M63").Value = -8486.4
$ws.Range("N63").Value = -5407

$ws.Range("H66").Value = 7986.846
$ws.Range("I66").Value = 9172.4
$ws.Range("J66").Value = 4035
$ws.Range("K66").Value = 45862
$ws.Range("L66").Value = 20175
$ws.Range("M66").Value = -42430
$ws.Range("N66").Value = -27039

$ws.Range("H101").Value = 31610
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 31610
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 31610
$ws.Range("N101").Value = -38100

$ws.Range("H102").Value = 1766.9062
$ws.Range("I102").Value = 1456.5
$ws.Range("J102").Value = 2284.25
$ws.Range("K102").Value = 1456.5
$ws.Range("L102").Value = 2284.25
$ws.Range("M102").Value = 165.5
$ws.Range("N102").Value = -5528.25

$ws.Range("H132").Value = 951.3099999999999
$ws.Range("I132").Value = 802.60254
$ws.Range("J132").Value = 1478.5454
$ws.Range("K132").Value = 2407.80762
$ws.Range("L132").Value = 4435.6362
$ws.Range("M132").Value = 122.19238
$ws.Range("N132").Value = -9495.636200000001

$ws.Range("H136").Value = 249098.11
$ws.Range("I136").Value = 212728.28
$ws.Range("J136").Value = 306408.12
$ws.Range("K136").Value = 638184.84
$ws.Range("L136").Value = 919224.36
$ws.Range("M136").Value = -635634.84
$ws.Range("N136").Value = -924324.36

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1529.125
$ws.Range("I107").Value = 1435.6538
$ws.Range("J107").Value = 1702.7142
$ws.Range("K107").Value = 1435.6538
$ws.Range("L107").Value = 1702.7142
$ws.Range("M107").Value = 484.3462
$ws.Range("N107").Value = -5542.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 609.48
$ws.Range("I16").Value = 573.3889
$ws.Range("J16").Value = 702.2857
$ws.Range("K16").Value = 573.3889
$ws.Range("L16").Value = 702.2857
$ws.Range("M16").Value = -286.3889
$ws.Range("N16").Value = -1276.2857

$ws.Range("H58").Value = 2212.679
$ws.Range("I58").Value = 2485.4905
$ws.Range("J58").Value = 1696.2858
$ws.Range("K58").Value = 2485.4905
$ws.Range("L58").Value = 1696.2858
$ws.Range("M58").Value = -2282.4905
$ws.Range("N58").Value = -2102.2858

$ws.Range("H105").Value = 575.2941
$ws.Range("I105").Value = 575.2941
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 575.2941
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1171.7059

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null

$ws.Range("H113").Value = 609.48
$ws.Range("I113").Value = 573.3889
$ws.Range("J113").Value = 702.2857
$ws.Range("K113").Value = 573.3889
$ws.Range("L113").Value = 702.2857
$ws.Range("M113").Value = 1596.6111
$ws.Range("N113").Value = -5042.2857

$ws.Range("H132").Value = 1611.44
$ws.Range("I132").Value = 933.5263
$ws.Range("J132").Value = 3758.1667
$ws.Range("K132").Value = 2800.5789
$ws.Range("L132").Value = 11274.5001
$ws.Range("M132").Value = -270.5789
$ws.Range("N132").Value = -16334.5001

$ws.Range("H136").Value = 2212.679
$ws.Range("I136").Value = 2485.4905
$ws.Range("J136").Value = 1696.2858
$ws.Range("K136").Value = 7456.4715
$ws.Range("L136").Value = 5088.857400000001
$ws.Range("M136").Value = -4906.4715
$ws.Range("N136").Value = -10188.8574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 2141.111
$ws.Range("I104").Value = 1403
$ws.Range("J104").Value = 2233.375
$ws.Range("K104").Value = 4209
$ws.Range("L104").Value = 6700.125
$ws.Range("M104").Value = -1588
$ws.Range("N104").Value = -11942.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 658.75
$ws.Range("I2").Value = 21.666666
$ws.Range("J2").Value = 1041
$ws.Range("K2").Value = 21.666666
$ws.Range("L2").Value = 1041
$ws.Range("M2").Value = 91.33333400000001
$ws.Range("N2").Value = -1267

$ws.Range("H132").Value = 3490.1428
$ws.Range("I132").Value = 3135.4333
$ws.Range("J132").Value = 4050.2104
$ws.Range("K132").Value = 9406.2999
$ws.Range("L132").Value = 12150.6312
$ws.Range("M132").Value = -6876.2999
$ws.Range("N132").Value = -17210.6312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1718
$ws.Range("I61").Value = 1677.7142
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1677.7142
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1475.7142
$ws.Range("N61").Value = -2404

$ws.Range("H93").Value = 1856.3529
$ws.Range("I93").Value = 1933.7778
$ws.Range("J93").Value = 1769.25
$ws.Range("K93").Value = 1933.7778
$ws.Range("L93").Value = 1769.25
$ws.Range("M93").Value = -685.7778000000001
$ws.Range("N93").Value = -4265.25

$ws.Range("H100").Value = 71433940
$ws.Range("I100").Value = 11800
$ws.Range("J100").Value = 111112900
$ws.Range("K100").Value = 11800
$ws.Range("L100").Value = 111112900
$ws.Range("M100").Value = -11259
$ws.Range("N100").Value = -111113982

$ws.Range("H113").Value = 1718
$ws.Range("I113").Value = 1677.7142
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1677.7142
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 492.2858000000001
$ws.Range("N113").Value = -6340

$ws.Range("H132").Value = 6352.3887
$ws.Range("I132").Value = 2656.9697
$ws.Range("J132").Value = 12159.477
$ws.Range("K132").Value = 7970.909100000001
$ws.Range("L132").Value = 36478.431
$ws.Range("M132").Value = -5440.909100000001
$ws.Range("N132").Value = -41538.431

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1942.678
$ws.Range("I132").Value = 1167.3103
$ws.Range("J132").Value = 2692.2
$ws.Range("K132").Value = 3501.9309
$ws.Range("L132").Value = 8076.599999999999
$ws.Range("M132").Value = -971.9309000000003
$ws.Range("N132").Value = -13136.6
